$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert new worksheet "margin_5_or_less" right before "gdp_andprezresults"
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("gdp_andprezresults")
$marginSheet = $wb.Worksheets.Add($refSheet)
$marginSheet.Name = "margin_5_or_less"

$marginHeaders = @("margin_flag", "stance", "n")
for ($col = 1; $col -le $marginHeaders.Length; $col++) {
    $marginSheet.Cells.Item(1, $col).Value = $marginHeaders[$col - 1]
}

$marginData = @(
    @("5_points_or_less", "not_sponsoring", 11),
    @("5_points_or_less", "sponsoring", 15),
    @("more_than_5_points", "not_sponsoring", 14),
    @("more_than_5_points", "sponsoring", 190),
    @("other", "sponsoring", 4)
)

for ($r = 0; $r -lt $marginData.Length; $r++) {
    $rowVals = $marginData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $marginSheet.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

$marginSheet.Range("A1:C1").Font.Bold = $true
$marginSheet.Range("A1:C1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) Append new worksheet "margin_5_or_less_withprez" at the very end
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$marginPrezSheet = $wb.Worksheets.Add($null, $lastSheet)
$marginPrezSheet.Name = "margin_5_or_less_withprez"

$marginPrezHeaders = @("p16winningparty", "margin_flag", "stance", "n")
for ($col = 1; $col -le $marginPrezHeaders.Length; $col++) {
    $marginPrezSheet.Cells.Item(1, $col).Value = $marginPrezHeaders[$col - 1]
}

$marginPrezData = @(
    @("D", "5_points_or_less", "not_sponsoring", 1),
    @("D", "5_points_or_less", "sponsoring", 11),
    @("D", "more_than_5_points", "not_sponsoring", 9),
    @("D", "more_than_5_points", "sponsoring", 178),
    @("D", "other", "sponsoring", 4),
    @("R", "5_points_or_less", "not_sponsoring", 10),
    @("R", "5_points_or_less", "sponsoring", 4),
    @("R", "more_than_5_points", "not_sponsoring", 5),
    @("R", "more_than_5_points", "sponsoring", 12)
)

for ($r = 0; $r -lt $marginPrezData.Length; $r++) {
    $rowVals = $marginPrezData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $marginPrezSheet.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

$marginPrezSheet.Range("A1:D1").Font.Bold = $true
$marginPrezSheet.Range("A1:D1").HorizontalAlignment = -4108
